# 2018-09-06 models relations issues in migrations
#
# - "product category" sheet gains a new "desc" column (C).
# - Selection/active-sheet bookkeeping is updated to match:
#     products sheet          -> selection moves to E1
#     stock operations sheet  -> selection moves to E2, no longer the active tab
#     product category sheet  -> selection moves to C3, becomes the active tab
#
# Order matters here: whichever sheet/range is selected LAST ends up marked as
# the active tab (tabSelected="1" on its sheetView + workbookView activeTab),
# so "product category" must be touched last.

$wb = $excel.ActiveWorkbook

$wsProducts = $wb.Worksheets.Item("products")
$wsCategory = $wb.Worksheets.Item("product category")
$wsStock    = $wb.Worksheets.Item("stock operations")

# products: just move the selection, no data changes.
$wsProducts.Range("E1").Select()

# stock operations: move the selection; this also drops its "active tab" status
# once product category is selected below.
$wsStock.Range("E2").Select()

# product category: add the new "desc" header in column C, then select C3 and
# make this the active sheet (must be the last selection of the script).
$wsCategory.Range("C1").Value = "desc"
$wsCategory.Range("C3").Select()
